$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

# Column G holds the "Recorded By" values (comma-separated list of recorders).
# Re-order each list so that any "System"/"system" entries come first,
# preserving the relative order of all other entries after them.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ($value -eq $null -or $value -eq "") {
        continue
    }

    $parts = $value -split ", "

    $systemParts = @()
    $otherParts = @()

    foreach ($part in $parts) {
        if ($part -eq "System") {
            $systemParts += $part
        } else {
            $otherParts += $part
        }
    }

    if ($systemParts.Count -gt 0) {
        $newValue = ($systemParts + $otherParts) -join ", "
        $cell.Value2 = $newValue
    }
}
